$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header for the "1h ahead CI" column to add a footnote marker ---
$ws.Range("C1").Value = "1h ahead CI*"

# --- Fill in the new LSTM row (row 6) results for "past p 168h, ws 1h forecast" ---
$ws.Range("C6").Value = "±2.98"
$ws.Range("D6").Value = "±4.24"
$ws.Range("E6").Value = "±4.93"

# --- Footnote about confidence intervals (italic) ---
$ws.Range("A9").Value = "* CI: 95% confidence intervals in MW"
$ws.Range("A9").Font.Italic = $true

# --- Section with ideas to improve the LSTM model (bold heading + bullet list) ---
$ws.Range("A11").Value = "Ideas to improve LSTM:"
$ws.Range("A11").Font.Bold = $true

$ws.Range("A12").Value = "Initialization of FFNN part"
$ws.Range("A13").Value = "Batch normalization"
$ws.Range("A14").Value = "Gradient clipping"

# --- Widen column C slightly to fit the new text ---
$ws.Columns.Item(3).ColumnWidth = 16

# --- Update the active selection, matching the saved workbook view ---
$ws.Range("E7").Select()
